$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'65.755.28"
$ws.Range("E2").Value = "  +0.36%  "
$ws.Range("D3").Value = "'2.674.87"
$ws.Range("E3").Value = "  +1.17%  "
$ws.Range("D5").Value = "'601.08"
$ws.Range("E5").Value = "  -0.36%  "
$ws.Range("D6").Value = "'157.25"
$ws.Range("E6").Value = "  +0.60%  "
$ws.Range("D7").Value = "'1.00"
$ws.Range("E7").Value = "  -0.05%  "
$ws.Range("D8").Value = "'0.613"
$ws.Range("E8").Value = "  +4.46%  "
$ws.Range("E9").Value = "  -0.14%  "
$ws.Range("E10").Value = "  +1.57%  "
$ws.Range("E11").Value = "  -0.17%  "
$ws.Range("E12").Value = "  -0.23%  "
$ws.Range("D13").Value = "'29.66"
$ws.Range("E13").Value = "  +0.05%  "
$ws.Range("D14").Value = "'0.0000198"
$ws.Range("E14").Value = "  +2.32%  "
$ws.Range("D15").Value = "'3.153.96"
$ws.Range("E15").Value = "  +1.12%  "
$ws.Range("D16").Value = "'65.550.02"
$ws.Range("E16").Value = "  +0.49%  "
$ws.Range("D17").Value = "'2.671.44"
$ws.Range("E17").Value = "  +0.47%  "
$ws.Range("D18").Value = "'12.63"
$ws.Range("E18").Value = "  +0.25%  "
$ws.Range("D19").Value = "'4.84"
$ws.Range("E19").Value = "  -0.34%  "
$ws.Range("D20").Value = "'7.61"
$ws.Range("E20").Value = "  +2.78%  "
$ws.Range("D21").Value = "'352.62"
$ws.Range("E21").Value = "  -1.16%  "
$ws.Range("D22").Value = "'1.00"
$ws.Range("E22").Value = "  -0.03%  "
$ws.Range("D23").Value = "'69.82"
$ws.Range("E23").Value = "  +0.44%  "
$ws.Range("D24").Value = "'0.0000112"
$ws.Range("E24").Value = "  +6.44%  "
$ws.Range("D25").Value = "'9.83"
$ws.Range("E25").Value = "  +4.98%  "
$ws.Range("E26").Value = "  -3.96%  "
$ws.Range("E27").Value = "  +2.05%  "
$ws.Range("E28").Value = "  -0.66%  "
$ws.Range("E29").Value = "  +1.00%  "
$ws.Range("D30").Value = "'544.10"
$ws.Range("E30").Value = "  +2.71%  "
$ws.Range("E31").Value = "  -0.08%  "
$ws.Range("E32").Value = "  -0.74%  "
$ws.Range("D33").Value = "'1.77"
$ws.Range("E33").Value = "  +0.15%  "
$ws.Range("D34").Value = "'6.62"
$ws.Range("E34").Value = "  +4.49%  "
$ws.Range("E35").Value = "  -0.77%  "
$ws.Range("E36").Value = "  -1.55%  "
$ws.Range("D37").Value = "'20.48"
$ws.Range("E37").Value = "  -0.49%  "
$ws.Range("E38").Value = "  +0.00%  "
$ws.Range("E39").Value = "  -0.35%  "
$ws.Range("D40").Value = "'157.96"
$ws.Range("E40").Value = "  -2.25%  "
$ws.Range("D41").Value = "'1.00"
$ws.Range("D42").Value = "'42.57"
$ws.Range("E42").Value = "  +1.40%  "
$ws.Range("D43").Value = "'165.78"
$ws.Range("E43").Value = "  +1.06%  "
$ws.Range("D44").Value = "'4.09"
$ws.Range("E44").Value = "  -1.00%  "
$ws.Range("D45").Value = "'0.0617"
$ws.Range("E45").Value = "  +1.87%  "
$ws.Range("E46").Value = "  -2.30%  "
$ws.Range("D47").Value = "'23.33"
$ws.Range("E47").Value = "  +2.50%  "
$ws.Range("E48").Value = "  -0.51%  "
$ws.Range("D49").Value = "'0.0261"
$ws.Range("E49").Value = "  -0.20%  "
$ws.Range("D50").Value = "'0.101"
$ws.Range("E50").Value = "  +2.58%  "
$ws.Range("D51").Value = "'20.06"
$ws.Range("E51").Value = "  +2.68%  "
